# Edit script: "New crime data collected" — update weekly CompStat report
# (9th Precinct) to the new reporting week (Volume 30, Number 20,
# covering 5/15/2023 - 5/21/2023) with refreshed crime-complaint figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / title updates -------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/15/2023  Through  5/21/2023"

# --- Crime Complaints table (rows 14-30) ------------------------------------
$ws.Range("C14").Value = "0"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = -100
$ws.Range("F14").Value = "0"
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = -100
$ws.Range("I14").Value = "0"
$ws.Range("J14").Value = 4
$ws.Range("K14").Value = -100
$ws.Range("L14").Value = -100
$ws.Range("M14").Value = -100
$ws.Range("N14").Value = -100
$ws.Range("C15").Value = "0"
$ws.Range("D15").Value = "0"
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 8
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = -20
$ws.Range("L15").Value = -20
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = -57.894736842105
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 250
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 46.666666666666
$ws.Range("I16").Value = 72
$ws.Range("J16").Value = 101
$ws.Range("K16").Value = -28.712871287128
$ws.Range("L16").Value = 44
$ws.Range("M16").Value = 18.032786885245
$ws.Range("N16").Value = -76.996805111821
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 91.666666666666
$ws.Range("I17").Value = 92
$ws.Range("J17").Value = 72
$ws.Range("K17").Value = 27.777777777777
$ws.Range("L17").Value = 50.819672131147
$ws.Range("M17").Value = 58.620689655172
$ws.Range("N17").Value = -54.679802955665
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 11
$ws.Range("E18").Value = -27.272727272727
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 31
$ws.Range("H18").Value = -35.483870967741
$ws.Range("I18").Value = 108
$ws.Range("J18").Value = 137
$ws.Range("K18").Value = -21.167883211678
$ws.Range("L18").Value = -14.285714285714
$ws.Range("M18").Value = 12.5
$ws.Range("N18").Value = -64.705882352941
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 24
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 79
$ws.Range("G19").Value = 68
$ws.Range("H19").Value = 16.176470588235
$ws.Range("I19").Value = 388
$ws.Range("J19").Value = 354
$ws.Range("K19").Value = 9.604519774011
$ws.Range("L19").Value = 77.168949771689
$ws.Range("M19").Value = 32.876712328767
$ws.Range("N19").Value = -29.19708029197
$ws.Range("C20").Value = "0"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 15
$ws.Range("J20").Value = 15
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = -40
$ws.Range("M20").Value = -11.764705882352
$ws.Range("N20").Value = -92.268041237113
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = -13.636363636363
$ws.Range("F21").Value = 147
$ws.Range("G21").Value = 132
$ws.Range("H21").Value = 11.363636363636
$ws.Range("I21").Value = 683
$ws.Range("J21").Value = 693
$ws.Range("K21").Value = -1.443001443001
$ws.Range("L21").Value = 38.539553752535
$ws.Range("M21").Value = 28.867924528301
$ws.Range("N21").Value = -56.962822936357
$ws.Range("C22").Value = "0"
$ws.Range("D22").Value = "0"
$ws.Range("E22").Value = "***.*"
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 7
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 75
$ws.Range("M22").Value = -12.5
$ws.Range("N22").Value = "***.*"
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 11
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 120
$ws.Range("I23").Value = 53
$ws.Range("J23").Value = 54
$ws.Range("K23").Value = -1.851851851851
$ws.Range("L23").Value = -29.333333333333
$ws.Range("M23").Value = 23.255813953488
$ws.Range("N23").Value = "***.*"
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 50
$ws.Range("E24").Value = -58
$ws.Range("F24").Value = 107
$ws.Range("G24").Value = 168
$ws.Range("H24").Value = -36.309523809523
$ws.Range("I24").Value = 557
$ws.Range("J24").Value = 836
$ws.Range("K24").Value = -33.373205741626
$ws.Range("L24").Value = 67.267267267267
$ws.Range("M24").Value = -6.070826306914
$ws.Range("N24").Value = "***.*"
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 17.142857142857
$ws.Range("I25").Value = 184
$ws.Range("J25").Value = 176
$ws.Range("K25").Value = 4.545454545454
$ws.Range("L25").Value = 73.584905660377
$ws.Range("M25").Value = 6.35838150289
$ws.Range("N25").Value = "***.*"
$ws.Range("C26").Value = "0"
$ws.Range("D26").Value = "0"
$ws.Range("E26").Value = "***.*"
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -66.666666666666
$ws.Range("I26").Value = 15
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = -16.666666666666
$ws.Range("L26").Value = 7.142857142857
$ws.Range("M26").Value = "***.*"
$ws.Range("N26").Value = "***.*"
$ws.Range("C27").Value = "0"
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -85.714285714285
$ws.Range("I27").Value = 21
$ws.Range("J27").Value = 34
$ws.Range("K27").Value = -38.235294117647
$ws.Range("L27").Value = -4.545454545454
$ws.Range("M27").Value = "***.*"
$ws.Range("N27").Value = "***.*"
$ws.Range("C28").Value = "0"
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 2
$ws.Range("J28").Value = 3
$ws.Range("K28").Value = -33.333333333333
$ws.Range("L28").Value = -33.333333333333
$ws.Range("M28").Value = 100
$ws.Range("N28").Value = -84.615384615384
$ws.Range("C29").Value = "0"
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = -100
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("I29").Value = 2
$ws.Range("J29").Value = 3
$ws.Range("K29").Value = -33.333333333333
$ws.Range("L29").Value = -33.333333333333
$ws.Range("M29").Value = 100
$ws.Range("N29").Value = -77.777777777777
$ws.Range("C30").Value = "0"
$ws.Range("D30").Value = "0"
$ws.Range("E30").Value = "***.*"
$ws.Range("F30").Value = "0"
$ws.Range("G30").Value = "0"
$ws.Range("H30").Value = "***.*"
$ws.Range("I30").Value = "0"
$ws.Range("J30").Value = 7
$ws.Range("K30").Value = -100
$ws.Range("L30").Value = -100
$ws.Range("M30").Value = "***.*"
$ws.Range("N30").Value = "***.*"

